$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the "4a88754c-d4d3-40bf-883e-05388a36cbfc.md" file.
# Status moves from "Handed back: in sync with en-US" to "Ready for handoff",
# and the Latest Handoff Date is refreshed.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-22 06:33:56"

# --- zh-cn sheet: row 3 is the same file. Status + Latest Handoff Datetime update.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-22 06:33:46"

# --- de-de sheet: row 3 is the same file. Status + Latest Handoff Datetime update.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-22 06:33:56"
